$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.983.17"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "2.666.98"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'522.96"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'144.32"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  +7.27%  "
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").Value = "'0.335"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D13").Value = "3.136.15"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "58.981.16"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "'21.02"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "2.686.32"
$ws.Range("E16").Value = "  -4.51%  "
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "'339.32"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").Value = "'10.37"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'64.40"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").Value = "'0.419"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "0.0₃0801"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "'6.67"
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'18.91"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "'150.60"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").Value = "'0.899"
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("D37").Value = "'0.872"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "'36.91"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").Value = "'0.617"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'275.61"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "'19.86"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "'0.0970"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").Value = "'10.66"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").Value = "2.052.83"
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("D49").Value = "'4.72"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "'18.83"
$ws.Range("E51").Value = "  -2.50%  "
